# The workbook's header row (row 1) used to start with an "id" column
# (A1="id", B1="nombre", ... M1="email"). This edit removes that leading
# "id" column entirely, so the remaining headers (nombre, descripcion,
# desc_sitio, horario, transporte, url, direccion, codpostal, latitud,
# longitud, telefono, email) shift one place to the left and the sheet's
# used range shrinks from A1:M1 to A1:L1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Delete()
